# Add a new "skos:definition" column (header in E7) and populate the
# definition text for the A320-neo row (E13), matching the style already
# used by the neighboring cells in those rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New header cell, same formatting as the rest of the header row (row 7).
$ws.Range("E7").Value = "skos:definition"

# New definition value for the A320-neo row, matching the formatting of
# its row neighbors (B13:D13).
$ws.Range("E13").Value = "An Airbus A320 that uses a new, more efficient engine model.  NEO stands for New Engine Option."
$ws.Range("D13").Copy()
$ws.Range("E13").PasteSpecial(-4122)  # xlPasteFormats
